$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registro de defectos")

# New test case CP08 (row 18)
$ws.Range("C18").Value = "CP08"
$ws.Range("D18").Value = "Creacion de objeto Carne y FrutaVerdura permite valores inferiores a 0 en cantidad y precio"
$ws.Range("F18").Value = "Grave"
$ws.Range("G18").Value = "Solicitado"
$ws.Range("H18").Value = "Crear restricciones en métodos set que no permitan valores inferiores a 0"

# New test case CP09 (row 19)
$ws.Range("C19").Value = "CP09"
$ws.Range("D19").Value = "Objeto Carne permite ingresar categorias distintas a A-B-C "
$ws.Range("F19").Value = "Grave"
$ws.Range("G19").Value = "Solicitado"
$ws.Range("H19").Value = "Crear restricciones que no permitan ingresar categorias fuera del modelo de negocios"

# Widen column D to fit the new, longer descriptions
$ws.Columns.Item(4).ColumnWidth = 91

# Move the active selection, as it ended up after the edits were made
$ws.Range("D25").Select()
